# Updated cryptos list on Sun Oct  8 05:41:00 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddress, $text) {
    # Force the cell to keep a literal text value (matching the workbook's
    # original inline-string cells) instead of Excel auto-coercing a
    # numeric-looking string (e.g. "211.83") into a float. The temporary
    # "@" (Text) number format prevents the coercion; resetting the style
    # back to "Normal" afterwards avoids leaving a stray style on the cell.
    $range = $ws.Range($rangeAddress)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "27.927.22"
$ws.Range("E2").Value = "  -0.02%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.630.89"
$ws.Range("E3").Value = "  -0.65%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.08%  "

# Row 5 - BNB
Set-TextValue "D5" "211.83"
$ws.Range("E5").Value = "  -0.82%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.17%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.04%  "

# Row 8 - Solana
Set-TextValue "D8" "23.35"
$ws.Range("E8").Value = "  -1.75%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -2.23%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.0613"
$ws.Range("E10").Value = "  -0.54%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.60%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue "D12" "1.863.23"
$ws.Range("E12").Value = "  -0.58%  "

# Row 13 - WrappedEther
Set-TextValue "D13" "1.634.08"
$ws.Range("E13").Value = "  -0.42%  "

# Row 14 - Polkadot
Set-TextValue "D14" "4.05"
$ws.Range("E14").Value = "  -0.91%  "

# Row 15 - Polygon
Set-TextValue "D15" "0.561"
$ws.Range("E15").Value = "  -2.23%  "

# Row 16 - Litecoin
Set-TextValue "D16" "65.64"
$ws.Range("E16").Value = "  -0.77%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "27.942.85"
$ws.Range("E17").Value = "  +0.07%  "

# Row 18 - BitcoinCash
Set-TextValue "D18" "230.67"
$ws.Range("E18").Value = "  -0.22%  "

# Row 19 - ShibaInu
Set-TextValue "D19" "0.0₃0723"
$ws.Range("E19").Value = "  -0.44%  "

# Row 20 - Chainlink
Set-TextValue "D20" "7.62"
$ws.Range("E20").Value = "  +0.02%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  +0.02%  "

# Row 22 - Avalanche
$ws.Range("E22").Value = "  -8.84%  "

# Row 23 - Uniswap
$ws.Range("E23").Value = "  -0.95%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -0.08%  "

# Row 25 - Monero
Set-TextValue "D25" "155.53"
$ws.Range("E25").Value = "  +2.23%  "

# Row 26 - Cosmos
Set-TextValue "D26" "6.92"
$ws.Range("E26").Value = "  -0.28%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  -0.50%  "

# Row 28 - EthereumClassic
Set-TextValue "D28" "15.57"
$ws.Range("E28").Value = "  -1.10%  "

# Row 29 - BinanceUSD
$ws.Range("E29").Value = "  +0.09%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -0.34%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  -0.93%  "

# Row 32 - Filecoin
Set-TextValue "D32" "3.39"
$ws.Range("E32").Value = "  +1.47%  "

# Row 33 - Maker
Set-TextValue "D33" "1.398.14"
$ws.Range("E33").Value = "  -1.74%  "

# Row 34 - InternetComputer(DFINITY)
Set-TextValue "D34" "3.06"
$ws.Range("E34").Value = "  -1.74%  "

# Row 35 - TrustWalletToken
$ws.Range("E35").Value = "  +13.63%  "

# Row 36 - LidoDAOToken
$ws.Range("E36").Value = "  -0.32%  "

# Row 37 - HuobiToken
Set-TextValue "D37" "2.37"
$ws.Range("E37").Value = "  +0.90%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  +1.99%  "

# Row 39 - ImmutableX
Set-TextValue "D39" "0.555"
$ws.Range("E39").Value = "  -0.60%  "

# Row 40 - ARBITRUM
Set-TextValue "D40" "0.863"
$ws.Range("E40").Value = "  -3.23%  "

# Row 41 - WEMIXToken
$ws.Range("E41").Value = "  -0.69%  "

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  +0.07%  "

# Rows 43-46 were reordered: MXToken moved from row 46 to row 43,
# pushing Aave, RenderToken and FraxShare down by one row, each with
# updated price/volume figures.
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D43" "2.29"
$ws.Range("E43").Value = "  +3.86%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D44" "66.54"
$ws.Range("E44").Value = "  -0.88%  "

$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D45" "1.82"
$ws.Range("E45").Value = "  -0.10%  "

$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D46" "5.45"
$ws.Range("E46").Value = "  +0.04%  "

# Row 47 - RocketPoolETH
Set-TextValue "D47" "1.776.98"
$ws.Range("E47").Value = "  -0.32%  "

# Row 48 - Quant
Set-TextValue "D48" "87.96"
$ws.Range("E48").Value = "  -0.93%  "

# Row 49 - BabyDogeCoin
Set-TextValue "D49" "0.0₆0103"
$ws.Range("E49").Value = "  -1.92%  "

# Row 50 - Algorand
Set-TextValue "D50" "0.0996"
$ws.Range("E50").Value = "  -1.16%  "

# Row 51 - Cronos
$ws.Range("E51").Value = "  -0.35%  "
